$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.ClearFormats()
}

Set-TextValue "D2" "47.123.79"
Set-TextValue "E2" "  +1.17%  "
Set-TextValue "D3" "2.487.80"
Set-TextValue "E3" "  +0.59%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "322.11"
Set-TextValue "E5" "  -0.21%  "
Set-TextValue "D6" "108.35"
Set-TextValue "E6" "  +2.20%  "
Set-TextValue "E7" "  +0.73%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.535"
Set-TextValue "E9" "  -0.98%  "
Set-TextValue "D10" "38.79"
Set-TextValue "E10" "  +7.06%  "
Set-TextValue "D11" "0.0812"
Set-TextValue "E11" "  -0.36%  "
Set-TextValue "E12" "  +0.33%  "
Set-TextValue "D13" "18.29"
Set-TextValue "E13" "  -0.80%  "
Set-TextValue "D14" "7.17"
Set-TextValue "E14" "  +0.75%  "
Set-TextValue "D15" "2.876.54"
Set-TextValue "E15" "  +0.33%  "
Set-TextValue "D16" "2.485.20"
Set-TextValue "E16" "  +2.45%  "
Set-TextValue "D17" "0.849"
Set-TextValue "E17" "  +0.31%  "
Set-TextValue "D18" "47.041.15"
Set-TextValue "E18" "  +1.16%  "
Set-TextValue "D19" "12.69"
Set-TextValue "E19" "  +0.02%  "
Set-TextValue "D20" "6.61"
Set-TextValue "E20" "  +1.96%  "
Set-TextValue "D21" "0.0₃0937"
Set-TextValue "E21" "  -0.25%  "
Set-TextValue "D22" "2.73"
Set-TextValue "E22" "  +14.32%  "
Set-TextValue "D23" "70.63"
Set-TextValue "E23" "  +0.00%  "
Set-TextValue "D24" "246.55"
Set-TextValue "E24" "  -0.84%  "
Set-TextValue "E25" "  +1.28%  "
Set-TextValue "E26" "  -0.04%  "
Set-TextValue "D27" "25.86"
Set-TextValue "E27" "  -1.40%  "
Set-TextValue "E28" "  +1.19%  "
Set-TextValue "D29" "10.02"
Set-TextValue "E29" "  +2.17%  "
Set-TextValue "D30" "0.141"
Set-TextValue "E30" "  +9.05%  "
Set-TextValue "D31" "35.05"
Set-TextValue "E31" "  +0.92%  "
Set-TextValue "D32" "49.94"
Set-TextValue "E32" "  +0.52%  "
Set-TextValue "D33" "19.98"
Set-TextValue "E33" "  +0.84%  "
Set-TextValue "D34" "5.40"
Set-TextValue "E34" "  +1.16%  "
Set-TextValue "D35" "0.0784"
Set-TextValue "E35" "  +2.00%  "
Set-TextValue "E36" "  +0.12%  "
Set-TextValue "E37" "  +2.34%  "
Set-TextValue "D38" "4.68"
Set-TextValue "D39" "2.96"
Set-TextValue "E39" "  +0.03%  "
Set-TextValue "E40" "  +0.30%  "
Set-TextValue "B41" "WEMIXToken"
Set-TextValue "C41" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D41" "2.23"
Set-TextValue "E41" "  -0.71%  "
Set-TextValue "B42" "Monero"
Set-TextValue "C42" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D42" "121.07"
Set-TextValue "E42" "  -2.54%  "
Set-TextValue "D43" "21.21"
Set-TextValue "E43" "  +0.95%  "
Set-TextValue "D44" "0.0296"
Set-TextValue "E44" "  +0.68%  "
Set-TextValue "D45" "1.988.98"
Set-TextValue "E45" "  +0.18%  "
Set-TextValue "D46" "3.04"
Set-TextValue "E46" "  +1.70%  "
Set-TextValue "E47" "  -2.66%  "
Set-TextValue "D48" "1.79"
Set-TextValue "E48" "  -2.70%  "
Set-TextValue "D49" "9.11"
Set-TextValue "E49" "  +0.23%  "
Set-TextValue "D50" "5.16"
Set-TextValue "E50" "  +0.19%  "
Set-TextValue "D51" "56.27"
Set-TextValue "E51" "  +2.34%  "
